$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 226. This shifts the
# existing rows 226-286 down to 227-287 (dimension becomes A1:R287) and
# copies the formatting (incl. the date style on column D) from the row
# above, matching Excel's native "Insert" behaviour.
$ws.Rows("226:226").Insert()

# Populate the newly inserted row 226 with the new weekly price record.
$ws.Cells.Item(226, 1).Value2  = 9
$ws.Cells.Item(226, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(226, 3).Value2  = "Metropolitana"
$ws.Cells.Item(226, 4).Value2  = 44855
$ws.Cells.Item(226, 5).Value2  = 13
$ws.Cells.Item(226, 6).Value2  = 100112026
$ws.Cells.Item(226, 7).Value2  = "Haba"
$ws.Cells.Item(226, 8).Value2  = "Sin especificar"
$ws.Cells.Item(226, 9).Value2  = "Primera"
$ws.Cells.Item(226, 10).Value2 = 150
$ws.Cells.Item(226, 11).Value2 = 8000
$ws.Cells.Item(226, 12).Value2 = 9000
$ws.Cells.Item(226, 13).Value2 = 8533
$ws.Cells.Item(226, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(226, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(226, 16).Value2 = 341
$ws.Cells.Item(226, 17).Value2 = 25
$ws.Cells.Item(226, 18).Value2 = "Hortaliza"
